$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for each data row (2-247).
# Every value of 45204 is being bumped to 45205 (one day later).
$ws.Range("C2:C247").Value = 45205
